$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new column K and add header "product_discount" ---
$ws.Columns.Item(11).Insert()
$ws.Range("K1").Value = "product_discount"
$ws.Range("K1").Font.Size = 12
$ws.Columns.Item(11).ColumnWidth = 15.5

# --- Column A: fill the previously empty "id" column with the same values as "product_code" (column B) ---
$ws.Range("A2").Value = 5001
$ws.Range("A3").Value = 5002
$ws.Range("A4").Value = 5003
$ws.Range("A5").Value = 5004
$ws.Range("A6").Value = 5005
$ws.Range("A7").Value = 5006
$ws.Range("A8").Value = 5007

# --- Column E: store product_price values as text instead of numbers ---
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "3045.45"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1768.18"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1795.45"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1754.55"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "4022.73"

# --- Highlight all data rows (A2:J8) in yellow ---
$ws.Range("A2:J8").Interior.Color = 65535

# --- Highlight the missing-price cells (E3, E4) in red instead ---
$ws.Range("E3").Interior.Color = 255
$ws.Range("E4").Interior.Color = 255
